# Updated after Python for engineers
#
# 1) The "datetimeFigureOut" date field cached in the slide master and every
#    slide layout moves on one day: 16/01/2023 -> 17/01/2023.
# 2) Slide 4 ("Doel van unit tests?") gets three of its bullet texts rewritten.

$p = $ppt.ActivePresentation

function Update-DateShape($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            $tf = $sh.TextFrame
            if ($tf.HasText) {
                $tr = $tf.TextRange
                if ($tr.Text -eq "16/01/2023") {
                    $tr.Text = "17/01/2023"
                }
            }
        }
    }
}

# --- Slide master date placeholder ---
$master = $p.SlideMaster
Update-DateShape $master.Shapes

# --- Every slide layout's date placeholder ---
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DateShape $layout.Shapes
}

# --- Slide 4 content bullet edits ---
$slide4 = $p.Slides.Item(4)
$contentShape = $slide4.Shapes.Item(2)
$bodyRange = $contentShape.TextFrame.TextRange

$p10 = $bodyRange.Paragraphs(10, 1)
$p10.Runs(1, 1).Text = "Testen is vastleggen."

$p12 = $bodyRange.Paragraphs(12, 1)
$p12.Runs(1, 1).Text = "Tests geven concrete voorbeelden van de functionaliteit."

$p13 = $bodyRange.Paragraphs(13, 1)
$p13.Runs(1, 1).Text = "Tests moedigen grondige inspectie aan: Wat verwacht ik? Wat kan geschrapt worden?"
